$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 trailing rows (17-19) that belonged to the removed "Neutrophils"
# sending-cluster block; this shifts the dimension from A1:T19 to A1:T16.
$ws.Range("A17:T19").Delete()

$data = @(
    @("ECs", "Oxt", "Avpr1a", "ECs", 1, 0.5, 0.1388385, 0.277677, 0.113840447003231, 0.0815056876720435, 2, 1, 2.5383755, 5.076751, 0.1536689326353429, 0.126906281076628, 0.35242424685675, 1.409696987427, 0.01749373998171682, 0.01034358370905221),
    @("ECs", "Oxt", "Avpr1a", "FAPs", 1, 0.5, 0.1388385, 0.277677, 0.113840447003231, 0.0815056876720435, 3, 1, 6.967000000000001, 20.901, 0.4217703226612587, 0.5224735624777741, 0.9672878295000003, 5.803726977000001, 0.04801452206445466, 0.04258456700021336),
    @("ECs", "Oxt", "Avpr1a", "MuSCs", 1, 0.5, 0.1388385, 0.277677, 0.113840447003231, 0.0815056876720435, 2, 1, 7.0130935, 14.026187, 0.4245607447033983, 0.3506201564455979, 0.97368738189975, 3.894749527599, 0.04833218495705949, 0.02857753696277793),
    @("FAPs", "Oxt", "Avpr1a", "ECs", 3, 1, 0.3578333333333334, 1.0735, 0.2934049749840486, 0.3151011992924826, 2, 1, 2.5383755, 5.076751, 0.1536689326353429, 0.126906281076628, 0.9083153664166668, 5.449892198500001, 0.04508722933569823, 0.03998832136499438),
    @("FAPs", "Oxt", "Avpr1a", "FAPs", 3, 1, 0.3578333333333334, 1.0735, 0.2934049749840486, 0.3151011992924826, 3, 1, 6.967000000000001, 20.901, 0.4217703226612587, 0.5224735624777741, 2.493024833333334, 22.43722350000001, 0.1237495109694407, 0.1646320461353625),
    @("FAPs", "Oxt", "Avpr1a", "MuSCs", 3, 1, 0.3578333333333334, 1.0735, 0.2934049749840486, 0.3151011992924826, 2, 1, 7.0130935, 14.026187, 0.4245607447033983, 0.3506201564455979, 2.509518624083334, 15.0571117445, 0.1245682346789096, 0.1104808317921258),
    @("Inflammatory-Mac", "Oxt", "Avpr1a", "ECs", 1, 0.3333333333333333, 0.04659966666666667, 0.139799, 0.03820933590851887, 0.04103477648802029, 2, 1, 2.5383755, 5.076751, 0.1536689326353429, 0.126906281076628, 0.1182874521748333, 0.709724713049, 0.005871587865767374, 0.005207570878905308),
    @("Inflammatory-Mac", "Oxt", "Avpr1a", "FAPs", 1, 0.3333333333333333, 0.04659966666666667, 0.139799, 0.03820933590851887, 0.04103477648802029, 3, 1, 6.967000000000001, 20.901, 0.4217703226612587, 0.5224735624777741, 0.3246598776666668, 2.921938899000001, 0.01611556393480842, 0.02143958585717516),
    @("Inflammatory-Mac", "Oxt", "Avpr1a", "MuSCs", 1, 0.3333333333333333, 0.04659966666666667, 0.139799, 0.03820933590851887, 0.04103477648802029, 2, 1, 7.0130935, 14.026187, 0.4245607447033983, 0.3506201564455979, 0.3268078194021667, 1.960846916413, 0.01622218410794307, 0.01438761975193982),
    @("MuSCs", "Oxt", "Avpr1a", "ECs", 1, 0.5, 0.113085, 0.22617, 0.0927238982656855, 0.0663869941723156, 2, 1, 2.5383755, 5.076751, 0.1536689326353429, 0.126906281076628, 0.2870521934175, 1.14820877367, 0.01424878247627601, 0.008424926542264347),
    @("MuSCs", "Oxt", "Avpr1a", "FAPs", 1, 0.5, 0.113085, 0.22617, 0.0927238982656855, 0.0663869941723156, 3, 1, 6.967000000000001, 20.901, 0.4217703226612587, 0.5224735624777741, 0.7878631950000002, 4.727179170000001, 0.03910818848992791, 0.03468544934740096),
    @("MuSCs", "Oxt", "Avpr1a", "MuSCs", 1, 0.5, 0.113085, 0.22617, 0.0927238982656855, 0.0663869941723156, 2, 1, 7.0130935, 14.026187, 0.4245607447033983, 0.3506201564455979, 0.7930756784475, 3.17230271379, 0.03936692729948157, 0.02327661828265029),
    @("Resolving-Mac", "Oxt", "Avpr1a", "ECs", 2, 0.6666666666666666, 0.5632320000000001, 1.689696, 0.461821343838516, 0.495971342375138, 2, 1, 2.5383755, 5.076751, 0.1536689326353429, 0.126906281076628, 1.429694309616, 8.578165857696, 0.07096759297588444, 0.06294187858141176),
    @("Resolving-Mac", "Oxt", "Avpr1a", "FAPs", 2, 0.6666666666666666, 0.5632320000000001, 1.689696, 0.461821343838516, 0.495971342375138, 3, 1, 6.967000000000001, 20.901, 0.4217703226612587, 0.5224735624777741, 3.924037344000001, 35.31633609600001, 0.194782537202627, 0.2591319141376222),
    @("Resolving-Mac", "Oxt", "Avpr1a", "MuSCs", 2, 0.6666666666666666, 0.5632320000000001, 1.689696, 0.461821343838516, 0.495971342375138, 2, 1, 7.0130935, 14.026187, 0.4245607447033983, 0.3506201564455979, 3.949998678192001, 23.699992069152, 0.1960712136600045, 0.1738975496561041),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value2 = $row[$j]
    }
}
